# Commit: "This is second change in ose branch"
# The workbook's single sheet (Sheet1) gets a literal text value written into
# A1. Excel will store this as a shared string (the diff adds
# xl/sharedStrings.xml with one <si> entry "This is second change" and sheet1
# gains a row with <c r="A1" t="s"><v>0</v></c> referencing it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "This is second change"
